$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resultat_etape")

# Insert a new row at row 11 (pushes ADR..VAL down by one), matching the
# banded formatting of the surrounding rows, and add the new code "ADMP"
# (-> "Admis") between "ADMI" (row 10) and "ADR" (row 11/now 12).
$ws.Rows.Item(10).Copy()
$ws.Rows.Item(11).Insert()

$ws.Cells.Item(11, 2).Value = "ADMP"
$ws.Cells.Item(11, 3).Value = "Admis"
